$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that should be bumped
# from 45188 (2023-09-19) to 45189 (2023-09-20) for every data row (2..369).
$lastRow = 369
$ws.Range("C2:C$lastRow").Value = 45189
